# Regenerate the per-UF ranking content for each analysis module.
# The underlying data refresh re-broke a few rank ties differently, so the
# "uf" column (column A) needs updating on every sheet to show the correct
# state abbreviation for its rank.

$wb = $excel.ActiveWorkbook

# --- qtd (sheet 1) ---
$ws = $wb.Worksheets.Item("qtd")
$ws.Range("A11").Value = "PA"
$ws.Range("A12").Value = "ES"
$ws.Range("A26").Value = "MT"
$ws.Range("A27").Value = "TO"

# --- tot-arrecad (sheet 2) ---
$ws = $wb.Worksheets.Item("tot-arrecad")
$ws.Range("A21").Value = "PI"
$ws.Range("A22").Value = "AL"
$ws.Range("A23").Value = "MT"
$ws.Range("A24").Value = "RN"
$ws.Range("A25").Value = "RO"
$ws.Range("A26").Value = "TO"
$ws.Range("A27").Value = "AP"

# --- avg-arrecad (sheet 3) ---
$ws = $wb.Worksheets.Item("avg-arrecad")
$ws.Range("A21").Value = "RO"
$ws.Range("A22").Value = "TO"
$ws.Range("A23").Value = "AP"
$ws.Range("A24").Value = "AL"
$ws.Range("A25").Value = "MT"
$ws.Range("A26").Value = "RN"
$ws.Range("A27").Value = "PI"

# --- max-arrecad (sheet 4) ---
$ws = $wb.Worksheets.Item("max-arrecad")
$ws.Range("A26").Value = "PI"
$ws.Range("A27").Value = "RN"

# --- tx-sucesso (sheet 5) ---
$ws = $wb.Worksheets.Item("tx-sucesso")
$ws.Range("A26").Value = "RN"
$ws.Range("A27").Value = "PI"
